# Update "Forecast Comparison" sheet with correct forecast output:
#  - insert a new "Week_Start_Date" column after "Week" (new column B)
#  - renumber week labels from zero-padded (W01) to unpadded (W1)
#  - shift existing columns (ASIN..is_holiday_week) right by one
#  - store is_holiday_week values as boolean cells

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B (ASIN), pushing ASIN..is_holiday_week to C..J
$ws.Columns.Item(2).Insert()

# New header for the inserted column
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Week start dates (Sunday) for weeks W1..W16, and un-padded week labels
$weekStarts = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

# Store Week_Start_Date values as text (not auto-converted Excel date serials)
$ws.Range("B2:B17").NumberFormat = "@"

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2

    # Un-pad week label (W01 -> W1)
    $ws.Cells.Item($row, 1).Value = "W" + ($i + 1)

    # New Week_Start_Date column (B)
    $ws.Cells.Item($row, 2).Value = $weekStarts[$i]

    # is_holiday_week (now column J) becomes a boolean value
    $ws.Cells.Item($row, 10).Value = $false
}
